# Add the "Cannabis expungement" easy-form entry to the sorted list.
#
# The source spreadsheet keeps a single sorted table (A2:B.. / name,url)
# and a new row for "Cannabis expungement" needs to be inserted in its
# correct alphabetically-sorted position (row 5), pushing every
# subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (shifts existing rows 5-29 down to 6-30)
$ws.Rows("5:5").Insert()

$name = "Cannabis expungement"
$url = "https://www.illinoislegalaid.org/legal-information/cannabis-expungement"

$ws.Range("A5").Value = $name
$ws.Range("B5").Value = $url

# Row insertion does not carry the pre-existing hyperlink anchors along
# with it, so rebuild the full hyperlink collection (cell -> URL, since
# every hyperlinked cell's text literally is its own target URL) in the
# correct, now-shifted, row order.
$ws.Hyperlinks.Delete()

$cells = @("B2","B12","B6","B4","B10","B27","B13","B28","B25","B7","B8","B9","B26","B30","B14","B16","B29","B22","B23","B21","B20","B18","B17","B5")

foreach ($c in $cells) {
    $ws.Hyperlinks.Add($ws.Range($c), $ws.Range($c).Value2)
}

# Hyperlinks.Add silently re-applies formatting through a fresh style
# record; put every touched cell explicitly back on the shared
# "Hyperlink" cell style so the visible formatting matches the rest of
# the column.
foreach ($c in $cells) {
    $ws.Range($c).Style = "Hyperlink"
}

# Restore the selection left behind by the editor.
$ws.Range("E10").Select()
